$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "67.405.18"  # D2: '67.412.94' -> '67.405.18'
$ws.Cells.Item(2, 5).Value = "  -1.21%  "  # E2: '  -1.23%  ' -> '  -1.21%  '
Set-TextValue $ws.Cells.Item(3, 4) "3.221.13"  # D3: '3.220.06' -> '3.221.13'
$ws.Cells.Item(3, 5).Value = "  -1.66%  "  # E3: '  -1.67%  ' -> '  -1.66%  '
$ws.Cells.Item(4, 5).Value = "  +0.00%  "  # E4: '  +0.03%  ' -> '  +0.00%  '
Set-TextValue $ws.Cells.Item(5, 4) "577.47"  # D5: '577.27' -> '577.47'
$ws.Cells.Item(5, 5).Value = "  -1.76%  "  # E5: '  -1.79%  ' -> '  -1.76%  '
Set-TextValue $ws.Cells.Item(6, 4) "182.77"  # D6: '182.52' -> '182.77'
$ws.Cells.Item(6, 5).Value = "  -1.68%  "  # E6: '  -1.36%  ' -> '  -1.68%  '
$ws.Cells.Item(7, 5).Value = "  +0.00%  "  # E7: '  +0.02%  ' -> '  +0.00%  '
Set-TextValue $ws.Cells.Item(8, 4) "0.603"  # D8: '0.601' -> '0.603'
$ws.Cells.Item(8, 5).Value = "  +0.19%  "  # E8: '  +0.09%  ' -> '  +0.19%  '
Set-TextValue $ws.Cells.Item(9, 4) "3.218.29"  # D9: '3.219.66' -> '3.218.29'
$ws.Cells.Item(9, 5).Value = "  -1.76%  "  # E9: '  -1.63%  ' -> '  -1.76%  '
$ws.Cells.Item(10, 5).Value = "  -3.50%  "  # E10: '  -3.58%  ' -> '  -3.50%  '
Set-TextValue $ws.Cells.Item(11, 4) "6.56"  # D11: '6.52' -> '6.56'
$ws.Cells.Item(11, 5).Value = "  -2.63%  "  # E11: '  -3.03%  ' -> '  -2.63%  '
Set-TextValue $ws.Cells.Item(12, 4) "0.412"  # D12: '0.411' -> '0.412'
$ws.Cells.Item(12, 5).Value = "  -1.57%  "  # E12: '  -1.54%  ' -> '  -1.57%  '
Set-TextValue $ws.Cells.Item(13, 4) "3.776.99"  # D13: '3.778.22' -> '3.776.99'
$ws.Cells.Item(13, 5).Value = "  -1.80%  "  # E13: '  -1.64%  ' -> '  -1.80%  '
$ws.Cells.Item(14, 5).Value = "  +0.02%  "  # E14: '  -0.05%  ' -> '  +0.02%  '
Set-TextValue $ws.Cells.Item(15, 4) "27.73"  # D15: '27.68' -> '27.73'
$ws.Cells.Item(15, 5).Value = "  -3.71%  "  # E15: '  -3.63%  ' -> '  -3.71%  '
Set-TextValue $ws.Cells.Item(16, 4) "67.472.30"  # D16: '67.474.76' -> '67.472.30'
$ws.Cells.Item(16, 5).Value = "  -1.13%  "  # E16: '  -1.14%  ' -> '  -1.13%  '
$ws.Cells.Item(17, 5).Value = "  -2.50%  "  # E17: '  -2.48%  ' -> '  -2.50%  '
Set-TextValue $ws.Cells.Item(18, 4) "3.245.17"  # D18: '3.222.87' -> '3.245.17'
$ws.Cells.Item(18, 5).Value = "  -0.94%  "  # E18: '  -1.52%  ' -> '  -0.94%  '
Set-TextValue $ws.Cells.Item(19, 4) "5.75"  # D19: '5.74' -> '5.75'
$ws.Cells.Item(19, 5).Value = "  -2.42%  "  # E19: '  -2.36%  ' -> '  -2.42%  '
Set-TextValue $ws.Cells.Item(20, 4) "13.40"  # D20: '13.39' -> '13.40'
$ws.Cells.Item(20, 5).Value = "  -2.01%  "  # E20: '  -1.85%  ' -> '  -2.01%  '
Set-TextValue $ws.Cells.Item(21, 4) "395.02"  # D21: '393.89' -> '395.02'
$ws.Cells.Item(21, 5).Value = "  +2.98%  "  # E21: '  +2.66%  ' -> '  +2.98%  '
Set-TextValue $ws.Cells.Item(22, 4) "7.55"  # D22: '7.54' -> '7.55'
Set-TextValue $ws.Cells.Item(23, 4) "1.00"  # D23: '0.999' -> '1.00'
$ws.Cells.Item(23, 5).Value = "  -0.05%  "  # E23: '  -0.09%  ' -> '  -0.05%  '
$ws.Cells.Item(24, 5).Value = "  -0.62%  "  # E24: '  -0.64%  ' -> '  -0.62%  '
$ws.Cells.Item(25, 5).Value = "  -0.82%  "  # E25: '  -0.73%  ' -> '  -0.82%  '
$ws.Cells.Item(26, 5).Value = "  -3.36%  "  # E26: '  -3.56%  ' -> '  -3.36%  '
$ws.Cells.Item(27, 5).Value = "  -0.28%  "  # E27: '  -0.21%  ' -> '  -0.28%  '
Set-TextValue $ws.Cells.Item(28, 4) "9.55"  # D28: '9.54' -> '9.55'
$ws.Cells.Item(28, 5).Value = "  -3.82%  "  # E28: '  -3.56%  ' -> '  -3.82%  '
$ws.Cells.Item(29, 5).Value = "  -0.02%  "  # E29: '  -0.03%  ' -> '  -0.02%  '
$ws.Cells.Item(30, 5).Value = "  -2.68%  "  # E30: '  -2.65%  ' -> '  -2.68%  '
Set-TextValue $ws.Cells.Item(31, 4) "5.56"  # D31: '5.55' -> '5.56'
$ws.Cells.Item(31, 5).Value = "  -4.37%  "  # E31: '  -4.08%  ' -> '  -4.37%  '
Set-TextValue $ws.Cells.Item(32, 4) "22.57"  # D32: '22.56' -> '22.57'
$ws.Cells.Item(32, 5).Value = "  -1.89%  "  # E32: '  -1.77%  ' -> '  -1.89%  '
Set-TextValue $ws.Cells.Item(33, 4) "6.94"  # D33: '6.93' -> '6.94'
$ws.Cells.Item(33, 5).Value = "  -4.77%  "  # E33: '  -4.57%  ' -> '  -4.77%  '
$ws.Cells.Item(34, 5).Value = "  +0.02%  "  # E34: '  +0.00%  ' -> '  +0.02%  '
Set-TextValue $ws.Cells.Item(35, 4) "1.25"  # D35: '1.24' -> '1.25'
$ws.Cells.Item(35, 5).Value = "  -3.21%  "  # E35: '  -3.01%  ' -> '  -3.21%  '
$ws.Cells.Item(36, 5).Value = "  -1.21%  "  # E36: '  -1.24%  ' -> '  -1.21%  '
$ws.Cells.Item(38, 5).Value = "  -0.27%  "  # E38: '  -0.16%  ' -> '  -0.27%  '
Set-TextValue $ws.Cells.Item(39, 4) "0.802"  # D39: '0.803' -> '0.802'
$ws.Cells.Item(39, 5).Value = "  -4.35%  "  # E39: '  -4.16%  ' -> '  -4.35%  '
Set-TextValue $ws.Cells.Item(40, 4) "26.19"  # D40: '26.15' -> '26.19'
$ws.Cells.Item(40, 5).Value = "  -2.42%  "  # E40: '  -2.50%  ' -> '  -2.42%  '
$ws.Cells.Item(41, 5).Value = "  -1.78%  "  # E41: '  -1.67%  ' -> '  -1.78%  '
$ws.Cells.Item(42, 5).Value = "  -4.36%  "  # E42: '  -4.40%  ' -> '  -4.36%  '
$ws.Cells.Item(43, 5).Value = "  -6.05%  "  # E43: '  -6.06%  ' -> '  -6.05%  '
Set-TextValue $ws.Cells.Item(44, 4) "0.0683"  # D44: '0.0682' -> '0.0683'
$ws.Cells.Item(44, 5).Value = "  -1.16%  "  # E44: '  -1.11%  ' -> '  -1.16%  '
$ws.Cells.Item(45, 5).Value = "  -2.55%  "  # E45: '  -2.49%  ' -> '  -2.55%  '
Set-TextValue $ws.Cells.Item(46, 4) "2.584.73"  # D46: '2.587.38' -> '2.584.73'
$ws.Cells.Item(46, 5).Value = "  -2.79%  "  # E46: '  -2.57%  ' -> '  -2.79%  '
$ws.Cells.Item(47, 2).Value = "Bittensor"  # B47: 'InjectiveProtocol' -> 'Bittensor'
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"  # C47: 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' -> 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Cells.Item(47, 4) "333.07"  # D47: '24.47' -> '333.07'
$ws.Cells.Item(47, 5).Value = "  -4.96%  "  # E47: '  -4.18%  ' -> '  -4.96%  '
$ws.Cells.Item(48, 2).Value = "InjectiveProtocol"  # B48: 'Bittensor' -> 'InjectiveProtocol'
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"  # C48: 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' -> 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Cells.Item(48, 4) "24.45"  # D48: '332.88' -> '24.45'
$ws.Cells.Item(48, 5).Value = "  -4.72%  "  # E48: '  -4.89%  ' -> '  -4.72%  '
$ws.Cells.Item(49, 5).Value = "  -3.08%  "  # E49: '  -2.97%  ' -> '  -3.08%  '
$ws.Cells.Item(50, 5).Value = "  -0.18%  "  # E50: '  -0.01%  ' -> '  -0.18%  '
$ws.Cells.Item(51, 5).Value = "  -2.06%  "  # E51: '  -2.09%  ' -> '  -2.06%  '
